$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New progress-log entries ---------------------------------------------
# Order matters: these are new shared strings, and the engine registers each
# unique string the first time it is written. The target workbook's
# sharedStrings table shows "2.5 on 09/07" (idx 17) created first, then
# "2 on 09/07" (idx 18), then "1 on 09/07" (idx 19) - so we write the cells
# in that same chronological order (G6 before H2 before F8; F10 reuses the
# already-registered "2.5 on 09/07" string).
$ws.Range("G6").Value  = "2.5 on 09/07"
$ws.Range("H2").Value  = "2 on 09/07"
$ws.Range("F8").Value  = "1 on 09/07"
$ws.Range("F10").Value = "2.5 on 09/07"

# --- Formatting sweep -------------------------------------------------------
# The whole C2:K12 block (minus the already-centered label columns, which
# this simply re-applies harmlessly) moves from "vertical-center only" to
# "horizontal+vertical center" to match the rest of the table.
$ws.Range("C2:K12").HorizontalAlignment = -4108
$ws.Range("C2:K12").VerticalAlignment = -4108

# Row 11 gains a B-column cell matching the style used by every other
# separator row (B3, B5, B7, B9).
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("B11").VerticalAlignment = -4107

# H1 becomes a (blank) centered cell, matching the new H column of data.
$ws.Range("H1").HorizontalAlignment = -4108

# --- Four new blank rows (13-16), styled like the other "B1"-style cells --
$ws.Range("B13:K16").HorizontalAlignment = -4108
$ws.Range("B13:K16").VerticalAlignment = -4107

# --- Column widths (G widens, new H column width defined) ------------------
$ws.Columns.Item(7).ColumnWidth = 16.7109375
$ws.Columns.Item(8).ColumnWidth = 16.42578125

# --- Selection moves to the new last cell, F16 ------------------------------
$ws.Range("F16").Select()
